# The 4 data rows (2-5) get cyclically rotated:
#   new row 2 <- old row 5
#   new row 3 <- old row 2
#   new row 4 <- old row 3
#   new row 5 <- old row 4
#
# Capture every row's raw values BEFORE any writes happen, then write the
# snapshots back in rotated order. Using Value2 (array) read/write both
# correctly preserves numbers/booleans/strings and clears cells that must
# end up blank in their destination row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = "A2:AY5"

$row2 = $ws.Range("A2:AY2").Value2
$row3 = $ws.Range("A3:AY3").Value2
$row4 = $ws.Range("A4:AY4").Value2
$row5 = $ws.Range("A5:AY5").Value2

$ws.Range("A2:AY2").Value2 = $row5
$ws.Range("A3:AY3").Value2 = $row2
$ws.Range("A4:AY4").Value2 = $row3
$ws.Range("A5:AY5").Value2 = $row4

# Columns Y and AA hold plain-text dates formatted like "2023-03-11". Excel's
# Value setter auto-converts such strings to date serials, so re-apply them
# as explicit text after the bulk rotation (ClearFormats removes the
# temporary "@" number format again so no stray style is left behind).
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue "Y2" "2023-03-11"
Set-TextValue "AA2" "2023-03-11"

Set-TextValue "Y3" "2010-09-03"
Set-TextValue "AA3" "2010-09-03"

Set-TextValue "Y4" "2023-03-11"
Set-TextValue "AA4" "2023-03-11"

Set-TextValue "Y5" "2023-03-11"
Set-TextValue "AA5" "2023-03-11"

Write-Output "rotation applied"
